$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Sheet1_2(19Nov)" to "Sheet1_2(20Nov)"
$trackSheet = $wb.Worksheets.Item("Sheet1_2(19Nov)")
$trackSheet.Name = "Sheet1_2(20Nov)"

# Fill in the L column (day 6 inventory values) for rows 8-13
$trackSheet.Range("L8").Value = 502
$trackSheet.Range("L9").Value = 74
$trackSheet.Range("L10").Value = 245
$trackSheet.Range("L11").Value = 354
$trackSheet.Range("L12").Value = 2
$trackSheet.Range("L13").Value = 34

# Update the view: scroll position and active cell/selection
$trackSheet.Activate()
$trackSheet.Application.ActiveWindow.ScrollColumn = 9
$trackSheet.Application.ActiveWindow.ScrollRow = 2
$trackSheet.Range("L13").Select()
